{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is `async (context) => { ... }`.\n//\n// The document's date header plus the 5x5 \"problems\" grid (25 division\n// problems, one every 4th table row) are updated in strict document\n// order. Several old values repeat (e.g. \"94\u00f76=\" appears twice) and some\n// new values collide with other old values (e.g. \"37\u00f75=\" is simultaneously\n// a source and a target), so replacements MUST be applied positionally\n// (one-to-one against the ordered list of non-empty paragraphs) rather\n// than via a global find/replace-all, which could mis-fire on the\n// duplicated / re-used text.\nconst replacements = [\n  [\"2024-05-19 Sunday\", \"2024-05-20 Monday\"],\n  [\"45\u00f74=\", \"96\u00f79=\"],\n  [\"17\u00f72=\", \"17\u00f75=\"],\n  [\"38\u00f75=\", \"89\u00f77=\"],\n  [\"54\u00f75=\", \"41\u00f77=\"],\n  [\"70\u00f72=\", \"53\u00f78=\"],\n  [\"11\u00f74=\", \"79\u00f76=\"],\n  [\"24\u00f75=\", \"37\u00f75=\"],\n  [\"73\u00f74=\", \"51\u00f73=\"],\n  [\"94\u00f76=\", \"80\u00f79=\"],\n  [\"28\u00f74=\", \"48\u00f78=\"],\n  [\"47\u00f78=\", \"70\u00f74=\"],\n  [\"37\u00f75=\", \"14\u00f73=\"],\n  [\"25\u00f72=\", \"68\u00f77=\"],\n  [\"39\u00f79=\", \"11\u00f76=\"],\n  [\"68\u00f73=\", \"20\u00f79=\"],\n  [\"43\u00f74=\", \"26\u00f75=\"],\n  [\"64\u00f73=\", \"24\u00f78=\"],\n  [\"35\u00f76=\", \"33\u00f76=\"],\n  [\"12\u00f74=\", \"15\u00f75=\"],\n  [\"91\u00f75=\", \"80\u00f77=\"],\n  [\"80\u00f77=\", \"35\u00f79=\"],\n  [\"94\u00f76=\", \"14\u00f79=\"],\n  [\"65\u00f79=\", \"50\u00f74=\"],\n  [\"31\u00f72=\", \"70\u00f79=\"],\n  [\"82\u00f73=\", \"82\u00f77=\"],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Keep only paragraphs that actually carry text (the grid has lots of\n// blank \"answer\" rows interleaved with the problem rows).\nconst targets = paragraphs.items.filter((p) => p.text.trim().length > 0);\n\nif (targets.length !== replacements.length) {\n  throw new Error(\n    `Expected ${replacements.length} non-empty paragraphs, found ${targets.length}`\n  );\n}\n\nfor (let i = 0; i < targets.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const para = targets[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      `Paragraph ${i}: expected \"${oldText}\" but found \"${para.text}\"`\n    );\n  }\n  // Replace the run's text in place (via the whole paragraph range) so\n  // existing run formatting (font, size, etc.) is preserved.\n  para.getRange().insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# $word.ActiveDocument is the open document.\n#\n# The document's date header plus the 5x5 \"problems\" grid (25 division\n# problems, one every 4th table row) are updated in strict document\n# order. Several old values repeat (e.g. \"94\u00f76=\" appears twice) and some\n# new values collide with other old values (e.g. \"37\u00f75=\" is simultaneously\n# a source and a target), so replacements MUST be applied positionally\n# (one-to-one against the ordered list of non-empty paragraphs) rather\n# than via a global Find/Replace, which could mis-fire on the\n# duplicated / re-used text.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-05-19 Sunday\", \"2024-05-20 Monday\"),\n    @(\"45\u00f74=\", \"96\u00f79=\"),\n    @(\"17\u00f72=\", \"17\u00f75=\"),\n    @(\"38\u00f75=\", \"89\u00f77=\"),\n    @(\"54\u00f75=\", \"41\u00f77=\"),\n    @(\"70\u00f72=\", \"53\u00f78=\"),\n    @(\"11\u00f74=\", \"79\u00f76=\"),\n    @(\"24\u00f75=\", \"37\u00f75=\"),\n    @(\"73\u00f74=\", \"51\u00f73=\"),\n    @(\"94\u00f76=\", \"80\u00f79=\"),\n    @(\"28\u00f74=\", \"48\u00f78=\"),\n    @(\"47\u00f78=\", \"70\u00f74=\"),\n    @(\"37\u00f75=\", \"14\u00f73=\"),\n    @(\"25\u00f72=\", \"68\u00f77=\"),\n    @(\"39\u00f79=\", \"11\u00f76=\"),\n    @(\"68\u00f73=\", \"20\u00f79=\"),\n    @(\"43\u00f74=\", \"26\u00f75=\"),\n    @(\"64\u00f73=\", \"24\u00f78=\"),\n    @(\"35\u00f76=\", \"33\u00f76=\"),\n    @(\"12\u00f74=\", \"15\u00f75=\"),\n    @(\"91\u00f75=\", \"80\u00f77=\"),\n    @(\"80\u00f77=\", \"35\u00f79=\"),\n    @(\"94\u00f76=\", \"14\u00f79=\"),\n    @(\"65\u00f79=\", \"50\u00f74=\"),\n    @(\"31\u00f72=\", \"70\u00f79=\"),\n    @(\"82\u00f73=\", \"82\u00f77=\")\n)\n\n$idx = 0\n$total = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $total; $i++) {\n    $p = $d.Paragraphs($i)\n    $r = $p.Range\n    # Paragraph text includes the trailing paragraph mark; strip it (and\n    # any stray whitespace) to compare against our expected values, and\n    # to detect the blank \"answer\" rows that must be left untouched.\n    $full = $r.Text\n    $current = $full.TrimEnd([char]13, [char]7)\n\n    if ($current.Length -eq 0) {\n        continue\n    }\n\n    if ($idx -ge $replacements.Length) {\n        continue\n    }\n\n    $pair = $replacements[$idx]\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    if ($current -ne $oldText) {\n        throw \"Paragraph $i`: expected '$oldText' but found '$current'\"\n    }\n\n    # Replace only the text portion, leaving the paragraph mark (and its\n    # formatting) untouched; run formatting on the remaining text is kept\n    # because we're writing into the same run range.\n    $textRange = $r.Duplicate\n    $textRange.MoveEnd(1, -1) | Out-Null\n    $textRange.Text = $newText\n\n    $idx++\n}\n\nif ($idx -ne $replacements.Length) {\n    throw \"Expected to apply $($replacements.Length) replacements but only applied $idx\"\n}\n"}
